# Add a new "LogMessage_NewFolder" row to the Constants sheet, and make
# that sheet the active tab with A2 selected (mirrors the author's manual
# verification of the existing folder sequence).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Constants")

$ws.Range("A2").Value = "LogMessage_NewFolder"
$ws.Range("B2").Value = "The folder does not exist. It was created automatically by the system."

# Match the "best fit" auto-sized widths Excel computed for the new columns
# (closest values reproducible through this host's column-width quantization).
$ws.Columns.Item(1).ColumnWidth = 21.5
$ws.Columns.Item(2).ColumnWidth = 61.333333333333336

$ws.Activate()
$ws.Range("A2").Select()
